$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values such as "30.074.74" or "1.008" as plain
# text (matching the source site formatting), never as numbers. Force the
# whole column range to a text format first so Excel does not auto-convert
# numeric-looking strings (e.g. "1.008", "12.50") into actual numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.074.74"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "2.103.20"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("D5").Value = "344.82"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "0.5180"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("D8").Value = "0.4430"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").Value = "0.09441"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").Value = "52.48"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").Value = "1.176"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "25.35"
$ws.Range("E12").Value = "  +3.97%  "
$ws.Range("D13").Value = "2.110.82"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "6.723"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "8.080"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "99.53"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "0.00001167"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "20.67"
$ws.Range("E19").Value = "  +6.09%  "
$ws.Range("D20").Value = "0.06704"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("E22").Value = "  -3.73%  "
$ws.Range("D23").Value = "30.167.21"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").Value = "12.73"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Value = "2.335"
$ws.Range("D26").Value = "2.359.44"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "22.11"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "164.09"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").Value = "2.540"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "134.07"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").Value = "1.156"
$ws.Range("E31").Value = "  -3.60%  "
$ws.Range("D32").Value = "0.1056"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").Value = "1.625"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "6.248"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "3.958"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "6.180"
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("D37").Value = "10.13"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").Value = "0.02575"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").Value = "0.06766"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Value = "0.2280"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "0.6956"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").Value = "12.50"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  +3.21%  "
$ws.Range("D44").Value = "0.6693"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("D45").Value = "14.16"
$ws.Range("E45").Value = "  -5.96%  "
$ws.Range("D46").Value = "2.280"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").Value = "3.640"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "0.00000000351"
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("D49").Value = "1.218"
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("D50").Value = "82.12"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "0.07185"
$ws.Range("E51").Value = "  -1.82%  "

# Restore the default (unstyled) cell style now that the text values are set,
# so no visible/number-format styling change is left behind on the cells.
$priceRange.Style = "Normal"
